# Fix the create_emp.xlsx test-data sheet so file names / expected
# behaviours line up correctly with each test case (datatype check for
# both inputs & response).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "create_emloyee_with_empty_data" should reference its own
# payload file and is a negative test case.
$ws.Range("B2").Value = "create_emloyee_with_empty_data.json"
$ws.Range("C2").Value = "negative"

# Row 3: strip the accidental leading space from the test name.
$ws.Range("A3").Value = "create_emloyee_with_existing_name"

# Row 5: "create_emloyee_with_valid_data" should reference its own
# payload file and is a positive test case.
$ws.Range("B5").Value = "create_emloyee_with_valid_data.json"
$ws.Range("C5").Value = "positive"
